# Apply the latest cryptos snapshot (Price / Volume(1h)) values.
# D-column prices that parse as plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the source inlineStr cells)
# instead of silently converting them to floating point numbers.
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range('D2').Value = '26.026.00'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '1.635.92'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').Value = "'0.992"
$ws.Range('E4').Value = '  -0.92%  '
$ws.Range('D5').Value = "'215.05"
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = "'0.994"
$ws.Range('E7').Value = '  -0.77%  '
$ws.Range('D8').Value = "'0.256"
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('D9').Value = "'0.0632"
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').Value = "'19.72"
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = "'4.24"
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').Value = '1.863.90'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '1.636.75'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = "'63.17"
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '25.999.10'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = "'192.39"
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('D23').Value = "'6.38"
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').Value = "'1.79"
$ws.Range('E25').Value = '  -1.77%  '
$ws.Range('D26').Value = "'141.49"
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = "'6.88"
$ws.Range('D29').Value = "'15.59"
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').Value = "'3.34"
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = "'0.906"
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('D37').Value = '1.143.65'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = "'2.48"
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('D41').Value = "'0.994"
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('D42').Value = "'5.58"
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = "'100.35"
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').Value = '1.774.81'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').Value = "'55.58"
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('D49').Value = "'1.46"
$ws.Range('E49').Value = '  +5.83%  '
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').Value = "'7.59"
$ws.Range('E51').Value = '  +0.45%  '
